$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
  @(1, -0.35844298829246668, 0.35739228422400515),
  @(2, -0.26350504993167334, 0.26044434774690117),
  @(3, -0.15749132295803392, 0.15656906413101623),
  @(4, -0.1445690642075661, 0.14374372332863139),
  @(5, -0.13774372361127796, 0.1360934167712653),
  @(6, -0.035121712964973462, 0.035105472264315374),
  @(7, -0.015105472609315385, 0.015094984881686813),
  @(8, 0.0049050147726985571, -0.0049165201877858777),
  @(9, 0.010916519896225552, -0.010938495908067836),
  @(10, 0.01693849561717542, -0.016943306007398462),
  @(11, 0.021443305722073092, -0.021462437058872297),
  @(12, 0.027462436768159293, -0.027586320345977189),
  @(13, 0.033586320057896302, -0.03365110726325149),
  @(14, 0.045651106953298992, -0.045777537362803322),
  @(15, 0.051777537077601465, -0.051969399306483588),
  @(16, 0.0072439643302946166, -0.0072530251436577764),
  @(17, 0.013253024861445972, -0.013265226269796315),
  @(18, -0.083887234575040992, 0.083783934343919952),
  @(19, -0.074783934623554149, 0.073992952641662058),
  @(20, -0.064992952928822234, 0.064821620127947455),
  @(21, -0.055821620416763196, 0.055590928186678212),
  @(22, -0.093953025076091379, 0.09363801115928716),
  @(23, -0.084638011449007955, 0.084127584307849723),
  @(24, -0.042127584729214185, 0.041999999576296965),
  @(25, -0.09497170409706257, 0.094722620805875124),
  @(26, -0.088722621097609533, 0.088402324878060767),
  @(27, -0.082402325171535118, 0.081308733238840691),
  @(28, -0.075308733539001693, 0.074550899324575148),
  @(29, -0.062550899652402236, 0.062174829053649461),
  @(30, -0.042174829414490578, 0.042020759030384713),
  @(31, -0.027020759375169234, 0.027000893938177839),
  @(32, -0.0060008943065730236, 0.0059999996891288987)
)

foreach ($row in $values) {
  $r = $row[0]
  $ws.Cells.Item($r, 1).Value = $row[1]
  $ws.Cells.Item($r, 2).Value = $row[2]
}

$ws.Columns.Item(1).ColumnWidth = 13.85
